$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 and Row 4 swap identity: AIG data moves to row 3, MetLife data moves to row 4.
# Company names (column B)
$ws.Range("B3").Value = "American International Group, I"
$ws.Range("B4").Value = "MetLife, Inc."

# Tickers (column C)
$ws.Range("C3").Value = "AIG"
$ws.Range("C4").Value = "MET"

# Row 2 (UnitedHealth / UNH) updated values
$ws.Range("D2").Value = 331.4
$ws.Range("E2").Value = 57.2
$ws.Range("F2").Value = 0.49
$ws.Range("N2").Value = 50.68470204858703

# Row 3 (now American International Group, I / AIG) updated values
$ws.Range("D3").Value = 77.16
$ws.Range("E3").Value = 45.7
$ws.Range("F3").Value = 1.32
$ws.Range("H3").Value = 50
$ws.Range("I3").Value = 50
$ws.Range("J3").Value = 46
$ws.Range("K3").Value = 53.2
$ws.Range("N3").Value = 50.68470204858703

# Row 4 (now MetLife, Inc. / MET) updated values
$ws.Range("D4").Value = 78.36
$ws.Range("E4").Value = 48.5
$ws.Range("F4").Value = 2.35
$ws.Range("H4").Value = 23
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 33
$ws.Range("K4").Value = 53.2
$ws.Range("N4").Value = 50.68470204858703

# Row 5 (Prudential Financial, Inc. / PRU) updated values
$ws.Range("D5").Value = 111.38
$ws.Range("E5").Value = 69.5
$ws.Range("F5").Value = 2.89
$ws.Range("H5").Value = 46
$ws.Range("I5").Value = 36
$ws.Range("K5").Value = 47.6
$ws.Range("N5").Value = 50.68470204858703
